$wb = $excel.ActiveWorkbook

# Rename the sheet "EnergyPlusV10" -> "Sheet1"
$ws = $wb.Worksheets.Item("EnergyPlusV10")
$ws.Name = "Sheet1"

# Remove the (unused) external reference / external link entirely.
$sources = $wb.LinkSources()
if ($sources) {
    foreach ($src in $sources) {
        $wb.BreakLink($src)
    }
}
